$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.05
$ws.Range("C5").Value = 0.01
$ws.Range("C6").Value = 0.01
$ws.Range("C7").Value = 0.1
$ws.Range("C9").Value = 0.05

$ws.Range("D9").Select()

Write-Output "Done"
